$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Data for rows 19-23: MSSV (col C) and Ho va ten (col D)
$data = @(
    @{ Row = 19; MSSV = "0712152"; Name = "Lê Long Hồ" },
    @{ Row = 20; MSSV = "0712163"; Name = "Võ Minh Hiển" },
    @{ Row = 21; MSSV = "0712174"; Name = "Nguyễn Văn Hiếu" },
    @{ Row = 22; MSSV = "0712178"; Name = "Nguyễn Ngọc Hoà" },
    @{ Row = 23; MSSV = "0712190"; Name = "Lê Gia Quốc Huy" }
)

foreach ($item in $data) {
    $r = $item.Row
    $ws.Cells.Item($r, 3).Value = $item.MSSV
    $ws.Cells.Item($r, 4).Value = $item.Name

    # Columns E..K (5..11) get a value of 1 formatted as a percentage
    $ws.Range($ws.Cells.Item($r, 5), $ws.Cells.Item($r, 11)).Value = 1
    $ws.Range($ws.Cells.Item($r, 5), $ws.Cells.Item($r, 11)).NumberFormat = "0%"
}

# Update the selected cell to match the committed selection
$ws.Range("K19").Select()
